$wb = $excel.ActiveWorkbook

# --- "Notes" sheet (first sheet) -------------------------------------------------
$notes = $wb.Worksheets.Item(1)

# Replace the old "Typically MS or PhD" note with a small list of valid values
# that also backs the new drop-down on the Data sheet.
$notes.Range("B5").Value = "MS"
$notes.Range("C5").Value = "PhD"
$notes.Range("D5").Value = "ME"

# Move the selection on the Notes sheet (it is no longer the active tab).
$notes.Range("E5").Select()

# --- "Data" sheet (second sheet) -------------------------------------------------
$data = $wb.Worksheets.Item(2)

# Give control over the date format used for the Start Date column.
$data.Range("C2").NumberFormat = "m/d/yy;@"

# Restrict Start Date entries to real dates (1/1/1900 .. 1/1/2100).
$data.Range("C2").Validation.Add(4, 1, 1, "1", "73051")

# Restrict the "Current Program" column to the values listed on the Notes sheet.
$data.Range("B1:B1048576").Validation.Add(3, 1, 1, '=Notes!$B$5:$D$5')

# Make "Data" the active sheet/tab and set its selection.
$data.Activate()
$data.Range("B8").Select()
